$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $style = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $style
}

# Row 2
Set-TextValue "D2" "26.379.96"
Set-TextValue "E2" "  -0.45%  "

# Row 3
Set-TextValue "D3" "1.723.29"
Set-TextValue "E3" "  -0.53%  "

# Row 4
Set-TextValue "D4" "0.9991"
Set-TextValue "E4" "  -0.08%  "

# Row 5
Set-TextValue "D5" "242.57"
Set-TextValue "E5" "  -1.84%  "

# Row 6
Set-TextValue "D6" "0.9998"
Set-TextValue "E6" "  -0.05%  "

# Row 7
Set-TextValue "D7" "0.4879"
Set-TextValue "E7" "  -0.26%  "

# Row 8
Set-TextValue "D8" "0.2591"
Set-TextValue "E8" "  -3.14%  "

# Row 9
Set-TextValue "D9" "0.06188"
Set-TextValue "E9" "  -0.53%  "

# Row 10
Set-TextValue "D10" "1.722.97"

# Row 11
Set-TextValue "D11" "0.06982"
Set-TextValue "E11" "  -1.21%  "

# Row 12
Set-TextValue "D12" "15.51"
Set-TextValue "E12" "  -1.08%  "

# Row 13
Set-TextValue "D13" "4.520"
Set-TextValue "E13" "  -2.90%  "

# Row 14
Set-TextValue "D14" "0.5972"
Set-TextValue "E14" "  -2.09%  "

# Row 15
Set-TextValue "D15" "77.17"
Set-TextValue "E15" "  -0.36%  "

# Row 16
Set-TextValue "E16" "  -0.10%  "

# Row 17
Set-TextValue "D17" "26.385.63"
Set-TextValue "E17" "  -0.39%  "

# Row 18
Set-TextValue "E18" "  -0.10%  "

# Row 19
Set-TextValue "D19" "0.000007193"
Set-TextValue "E19" "  +0.30%  "

# Row 20
Set-TextValue "E20" "  -1.68%  "

# Row 21
Set-TextValue "D21" "1.945.02"
Set-TextValue "E21" "  -0.72%  "

# Row 22
Set-TextValue "D22" "4.439"
Set-TextValue "E22" "  -1.94%  "

# Row 23
Set-TextValue "D23" "8.499"
Set-TextValue "E23" "  -3.26%  "

# Row 24
Set-TextValue "D24" "5.087"
Set-TextValue "E24" "  -3.28%  "

# Row 25
Set-TextValue "D25" "138.14"
Set-TextValue "E25" "  -0.99%  "

# Row 26
Set-TextValue "E26" "  -1.30%  "

# Row 27
Set-TextValue "E27" "  -0.23%  "

# Row 28
Set-TextValue "D28" "106.34"
Set-TextValue "E28" "  -1.56%  "

# Row 29
Set-TextValue "D29" "1.724"
Set-TextValue "E29" "  -2.91%  "

# Row 30
Set-TextValue "E30" "  -1.75%  "

# Row 31
Set-TextValue "E31" "  -0.40%  "

# Row 32
Set-TextValue "D32" "3.652"
Set-TextValue "E32" "  -1.25%  "

# Row 33
Set-TextValue "D33" "0.04499"
Set-TextValue "E33" "  -1.86%  "

# Row 34
Set-TextValue "D34" "2.605"
Set-TextValue "E34" "  -0.45%  "

# Row 35
Set-TextValue "D35" "0.9972"
Set-TextValue "E35" "  -0.96%  "

# Row 36
Set-TextValue "D36" "0.6234"
Set-TextValue "E36" "  -2.23%  "

# Row 37
Set-TextValue "D37" "0.9304"
Set-TextValue "E37" "  +3.36%  "

# Row 38
Set-TextValue "D38" "1.961"
Set-TextValue "E38" "  -2.92%  "

# Row 39
Set-TextValue "D39" "2.382"
Set-TextValue "E39" "  -0.70%  "

# Row 40
Set-TextValue "D40" "0.9988"
Set-TextValue "E40" "  -0.52%  "

# Row 41
Set-TextValue "B41" "VeChain"
Set-TextValue "C41" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.01474"
Set-TextValue "E41" "  -2.28%  "

# Row 42
Set-TextValue "B42" "Quant"
Set-TextValue "C42" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D42" "100.63"
Set-TextValue "E42" "  -0.68%  "

# Row 43
Set-TextValue "D43" "5.457"
Set-TextValue "E43" "  +0.09%  "

# Row 44
Set-TextValue "D44" "0.3843"
Set-TextValue "E44" "  -1.17%  "

# Row 45
Set-TextValue "D45" "6.877"
Set-TextValue "E45" "  -1.06%  "

# Row 46
Set-TextValue "D46" "0.1164"
Set-TextValue "E46" "  -1.70%  "

# Row 47
Set-TextValue "D47" "0.05369"
Set-TextValue "E47" "  -0.30%  "

# Row 48
Set-TextValue "D48" "30.19"
Set-TextValue "E48" "  -1.23%  "

# Row 49
Set-TextValue "D49" "7.680"
Set-TextValue "E49" "  -1.53%  "

# Row 50
Set-TextValue "B50" "Aave"
Set-TextValue "C50" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "51.01"
Set-TextValue "E50" "  -0.63%  "

# Row 51
Set-TextValue "B51" "NEARProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "1.220"
Set-TextValue "E51" "  -2.32%  "
